$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text in the source data even when it looks numeric
# (e.g. "37.827.35", "170.70"). Excel auto-converts plain numeric-looking
# strings assigned via Range.Value into real numbers (dropping formatting like
# trailing zeros), so force literal text with a leading apostrophe and then
# reset the cell style back to Normal (the apostrophe trick flags the cell as
# quote-prefixed, which would otherwise leave a stray style behind).

$ws.Range('D2').Value = "'37.827.35"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.24%  '
$ws.Range('D3').Value = "'2.088.51"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'234.36"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = "'58.98"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.75%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'0.392"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').Value = '  +2.34%  '
$ws.Range('E11').Value = '  +3.25%  '
$ws.Range('D12').Value = "'2.395.48"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = "'14.76"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('D14').Value = "'21.39"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.91%  '
$ws.Range('D15').Value = "'0.771"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.57%  '
$ws.Range('D16').Value = "'5.31"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = "'2.085.98"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('D18').Value = "'37.735.10"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('D19').Value = "'6.24"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.87%  '
$ws.Range('D20').Value = "'71.61"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.82%  '
$ws.Range('D21').Value = "'0.0₃0831"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.52%  '
$ws.Range('D22').Value = "'228.59"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = "'2.41"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('D26').Value = "'170.70"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.12%  '
$ws.Range('D28').Value = "'9.06"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.17%  '
$ws.Range('D29').Value = "'1.42"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').Value = "'19.58"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.39%  '
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('E32').Value = '  +3.90%  '
$ws.Range('E33').Value = '  +2.49%  '
$ws.Range('E34').Value = '  +3.62%  '
$ws.Range('D35').Value = "'2.51"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').Value = "'3.50"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.40%  '
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').Value = "'5.46"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.35%  '
$ws.Range('D40').Value = "'0.0986"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.32%  '
$ws.Range('D41').Value = "'99.57"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.66%  '
$ws.Range('D42').Value = "'2.94"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').Value = "'0.0216"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('D44').Value = "'1.462.25"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('E45').Value = '  +1.21%  '
$ws.Range('D46').Value = "'4.30"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.71%  '
$ws.Range('D47').Value = "'16.16"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.55%  '
$ws.Range('E48').Value = '  +5.29%  '
$ws.Range('D49').Value = "'7.49"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.61%  '
$ws.Range('E50').Value = '  +2.80%  '
$ws.Range('D51').Value = "'47.62"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.20%  '
